$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 223 (existing rows 223:241 shift down to 226:244).
$ws.Rows.Item(223).Insert()
$ws.Rows.Item(223).Insert()
$ws.Rows.Item(223).Insert()

# New row 223: Kurakata / Especial
$ws.Range("A223").Value = 11
$ws.Range("B223").Value = "Vega Monumental Concepción"
$ws.Range("C223").Value = "Bíobío"
$ws.Range("D223").Value = 44931
$ws.Range("E223").Value = 8
$ws.Range("F223").Value = "Fruta"
$ws.Range("G223").Value = 100103
$ws.Range("H223").Value = "Frutos de hueso (carozo)"
$ws.Range("I223").Value = 100103004
$ws.Range("J223").Value = "Durazno"
$ws.Range("K223").Value = "Kurakata"
$ws.Range("L223").Value = "Especial"
$ws.Range("M223").Value = 50
$ws.Range("N223").Value = 17000
$ws.Range("O223").Value = 17000
$ws.Range("P223").Value = 17000
$ws.Range("Q223").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R223").Value = "Región de O'Higgins"
$ws.Range("S223").Value = 1133
$ws.Range("T223").Value = 15

# New row 224: Kurakata / Primera
$ws.Range("A224").Value = 11
$ws.Range("B224").Value = "Vega Monumental Concepción"
$ws.Range("C224").Value = "Bíobío"
$ws.Range("D224").Value = 44931
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = "Fruta"
$ws.Range("G224").Value = 100103
$ws.Range("H224").Value = "Frutos de hueso (carozo)"
$ws.Range("I224").Value = 100103004
$ws.Range("J224").Value = "Durazno"
$ws.Range("K224").Value = "Kurakata"
$ws.Range("L224").Value = "Primera"
$ws.Range("M224").Value = 50
$ws.Range("N224").Value = 15000
$ws.Range("O224").Value = 15000
$ws.Range("P224").Value = 15000
$ws.Range("Q224").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R224").Value = "Región de O'Higgins"
$ws.Range("S224").Value = 1000
$ws.Range("T224").Value = 15

# New row 225: Kurakata / Segunda
$ws.Range("A225").Value = 11
$ws.Range("B225").Value = "Vega Monumental Concepción"
$ws.Range("C225").Value = "Bíobío"
$ws.Range("D225").Value = 44931
$ws.Range("E225").Value = 8
$ws.Range("F225").Value = "Fruta"
$ws.Range("G225").Value = 100103
$ws.Range("H225").Value = "Frutos de hueso (carozo)"
$ws.Range("I225").Value = 100103004
$ws.Range("J225").Value = "Durazno"
$ws.Range("K225").Value = "Kurakata"
$ws.Range("L225").Value = "Segunda"
$ws.Range("M225").Value = 50
$ws.Range("N225").Value = 13000
$ws.Range("O225").Value = 13000
$ws.Range("P225").Value = 13000
$ws.Range("Q225").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R225").Value = "Región de O'Higgins"
$ws.Range("S225").Value = 867
$ws.Range("T225").Value = 15
